$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): Mean / Variance pairs in columns H:I and K:L
$ws.Range("H1").Value = "Mean"
$ws.Range("I1").Value = "Variance"
$ws.Range("K1").Value = "Mean"
$ws.Range("L1").Value = "Variance"

# Data values for columns H (Mean-like) and I (Variance-like)
$hValues = @(
    -0.020049145985394701,
    -0.55407472075987496,
    -1.7582004540570799,
    -2.7030293227134998,
    -3.9636587916359201,
    -5.45401021605125
)
$iValues = @(
    0.00154883344613096,
    0.00347865694172935,
    0.011903179497765501,
    0.0274058370014344,
    0.035884466763871699,
    0.052305454467623198
)
$kValues = @(
    -0.0219183048009872,
    0.68309141956960595,
    1.56758007797856,
    2.3124645283159699,
    3.18175477274134,
    4.3161427320140797
)
$lValues = @(
    0.0014086166539250601,
    0.0013433216525311001,
    0.0015412683601421101,
    0.0015530561628733499,
    0.00156812179785935,
    0.00223336432352688
)

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $hValues[$i]
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 11).Value = $kValues[$i]
    $ws.Cells.Item($row, 12).Value = $lValues[$i]
}

# Apply number formats matching the existing "Variance" column (B) style for H and K
$ws.Range("H2:H7").NumberFormat = "0.000"
$ws.Range("K2:K7").NumberFormat = "0.000"

# Apply new number format for I and L columns
$ws.Range("I2:I7").NumberFormat = "0.0000"
$ws.Range("L2:L7").NumberFormat = "0.0000"

# Update selection to match the new target range
$ws.Range("K1:L7").Select()
